$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold values such as "1.002" or "0.06786" that
# Excel would otherwise auto-convert to plain numbers (dropping the
# original textual formatting). Force those specific cells to Text
# format before writing, so the literal digit string is preserved,
# matching the inline string values from the source data feed.

foreach ($addr in @("D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D18", "D19", "D20", "D21", "D23", "D26", "D27", "D29", "D30", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D45", "D46", "D47", "D48", "D49", "D50", "D51")) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = "23.367.14"
$ws.Range("E2").Value = "  -1.54%  "

# Row 3
$ws.Range("D3").Value = "1.633.51"
$ws.Range("E3").Value = "  -1.25%  "

# Row 4
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.06%  "

# Row 5
$ws.Range("D5").Value = "1.002"
$ws.Range("E5").Value = "  +0.02%  "

# Row 6
$ws.Range("D6").Value = "299.55"
$ws.Range("E6").Value = "  -1.29%  "

# Row 7
$ws.Range("D7").Value = "0.3775"
$ws.Range("E7").Value = "  -0.75%  "

# Row 8
$ws.Range("D8").Value = "50.21"
$ws.Range("E8").Value = "  -1.65%  "

# Row 9
$ws.Range("D9").Value = "0.3517"
$ws.Range("E9").Value = "  -2.88%  "

# Row 10
$ws.Range("D10").Value = "0.08045"
$ws.Range("E10").Value = "  -2.02%  "

# Row 11
$ws.Range("D11").Value = "1.204"
$ws.Range("E11").Value = "  -3.52%  "

# Row 12
$ws.Range("D12").Value = "1.002"
$ws.Range("E12").Value = "  +0.07%  "

# Row 13
$ws.Range("D13").Value = "21.91"
$ws.Range("E13").Value = "  -3.43%  "

# Row 14
$ws.Range("D14").Value = "6.316"
$ws.Range("E14").Value = "  -3.16%  "

# Row 15
$ws.Range("D15").Value = "7.241"
$ws.Range("E15").Value = "  -2.80%  "

# Row 16
$ws.Range("D16").Value = "0.00001201"
$ws.Range("E16").Value = "  -2.84%  "

# Row 17
$ws.Range("D17").Value = "1.635.47"
$ws.Range("E17").Value = "  -1.03%  "

# Row 18
$ws.Range("D18").Value = "95.82"
$ws.Range("E18").Value = "  -1.58%  "

# Row 19
$ws.Range("D19").Value = "0.06959"
$ws.Range("E19").Value = "  -0.94%  "

# Row 20
$ws.Range("D20").Value = "6.694"
$ws.Range("E20").Value = "  -1.49%  "

# Row 21
$ws.Range("D21").Value = "17.28"
$ws.Range("E21").Value = "  -2.39%  "

# Row 22
$ws.Range("E22").Value = "  -0.01%  "

# Row 23
$ws.Range("D23").Value = "12.27"
$ws.Range("E23").Value = "  -4.68%  "

# Row 24
$ws.Range("D24").Value = "23.383.90"
$ws.Range("E24").Value = "  -1.46%  "

# Row 25
$ws.Range("E25").Value = "  -2.24%  "

# Row 26
$ws.Range("D26").Value = "2.876"
$ws.Range("E26").Value = "  -5.77%  "

# Row 27
$ws.Range("D27").Value = "20.77"
$ws.Range("E27").Value = "  -2.34%  "

# Row 28
$ws.Range("E28").Value = "  +0.19%  "

# Row 29
$ws.Range("D29").Value = "5.184"
$ws.Range("E29").Value = "  -0.86%  "

# Row 30
$ws.Range("D30").Value = "132.09"
$ws.Range("E30").Value = "  -1.71%  "

# Row 31
$ws.Range("D31").Value = "1.819.58"
$ws.Range("E31").Value = "  -1.01%  "

# Row 32
$ws.Range("D32").Value = "6.807"
$ws.Range("E32").Value = "  -1.80%  "

# Row 33
$ws.Range("D33").Value = "2.128"
$ws.Range("E33").Value = "  -2.96%  "

# Row 34
$ws.Range("D34").Value = "11.31"
$ws.Range("E34").Value = "  -3.84%  "

# Row 35
$ws.Range("D35").Value = "0.9708"
$ws.Range("E35").Value = "  -9.55%  "

# Row 36
$ws.Range("D36").Value = "0.02697"
$ws.Range("E36").Value = "  -4.08%  "

# Row 37
$ws.Range("D37").Value = "0.08717"
$ws.Range("E37").Value = "  -1.11%  "

# Row 38
$ws.Range("D38").Value = "0.2416"
$ws.Range("E38").Value = "  -3.90%  "

# Row 39
$ws.Range("D39").Value = "5.859"
$ws.Range("E39").Value = "  -3.94%  "

# Row 40
$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").Value = "0.06786"
$ws.Range("E40").Value = "  -4.29%  "

# Row 41
$ws.Range("B41").Value = "Aptos"
$ws.Range("C41").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D41").Value = "12.91"
$ws.Range("E41").Value = "  -0.66%  "

# Row 42
$ws.Range("D42").Value = "0.6803"
$ws.Range("E42").Value = "  -3.28%  "

# Row 43
$ws.Range("D43").Value = "1.299"
$ws.Range("E43").Value = "  -2.87%  "

# Row 44
$ws.Range("E44").Value = "  -1.80%  "

# Row 45
$ws.Range("D45").Value = "1.001"
$ws.Range("E45").Value = "  +0.04%  "

# Row 46
$ws.Range("D46").Value = "0.6295"
$ws.Range("E46").Value = "  -3.31%  "

# Row 47
$ws.Range("D47").Value = "2.230"
$ws.Range("E47").Value = "  -3.88%  "

# Row 48
$ws.Range("D48").Value = "3.897"
$ws.Range("E48").Value = "  -1.60%  "

# Row 49
$ws.Range("D49").Value = "0.07675"
$ws.Range("E49").Value = "  -3.45%  "

# Row 50
$ws.Range("D50").Value = "126.69"
$ws.Range("E50").Value = "  -1.25%  "

# Row 51
$ws.Range("D51").Value = "1.204"
$ws.Range("E51").Value = "  +0.88%  "
